# Apply the "ready to change generation of blocks and javaScript in code.js" edit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Wipe the existing used range completely (content + per-row/per-cell formatting)
# so we can rebuild the new layout from a clean sheet.
$ws.Range("A1:C28").EntireRow.Delete() | Out-Null

# --- Column A (day markers) ---
$ws.Range("A1").Value = "Saturday"
$ws.Range("A2").Value = "done"
$ws.Range("A3").Value = "done"
$ws.Range("A4").Value = "Sunday-Monday"
$ws.Range("A5").Value = "done"
$ws.Range("A13").Value = "After Monday"

# --- Column B (task notes) ---
$ws.Range("B2").Value = "remove rename Type and Method;"
$ws.Range("B3").Value = "wire up client methods to update Type and Method;"
$ws.Range("B5").Value = "write code methods listed below"
$ws.Range("B6").Value = "change method XML and JavaScript generators to handle Method type and params;"
$ws.Range("B8").Value = "change type XML and JavaScript generators to handle base types."
$ws.Range("B10").Value = "test TI rename"
$ws.Range("B11").Value = "test"
$ws.Range("B12").Value = "Check out (in code.js) self.renameTypeInActiveComic and self.renameMethodInActiveType to be sure everything's being done."
$ws.Range("B14").Value = "Save project to DB"

# --- Column C (extra detail notes) ---
$ws.Range("C7").Value = "in m_functionGenerateBlocksMethodFunctionString and m_functionGenerateJavaScriptMethodFunctionString"
$ws.Range("C9").Value = "in m_functionGenerateBlocksTypeNewFunctionString and m_functionGenerateJavaScriptTypeNewFunctionString"

# --- "methods in code.js to test" block (rows 16-22) ---
$ws.Range("B16").Value = "methods in code.js to test"
$ws.Range("B16").Font.Bold = $true
$ws.Range("B16").Font.Underline = $true

$ws.Range("B17").Value = "self.removeEvent"
$ws.Range("B17").Font.Strikethrough = $true

$ws.Range("B18").Value = "self.renameEvent"
$ws.Range("B18").Font.Strikethrough = $true

$ws.Range("B19").Value = "m_functionRemove_Type_Event"
$ws.Range("B19").Font.Strikethrough = $true

$ws.Range("B20").Value = "self.replaceMethod"
$ws.Range("B20").Font.Strikethrough = $true

$ws.Range("B21").Value = "self.replaceType"
$ws.Range("B21").Font.Strikethrough = $true

$ws.Range("B22").Value = "self.replaceProperty"
$ws.Range("B22").Font.Strikethrough = $true

# Column A best-fits its (now much shorter) contents.
$ws.Columns.Item(1).AutoFit() | Out-Null

# Selection ends on B29, matching the saved view state.
$ws.Range("B29").Select() | Out-Null
